$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 5136.6343
$ws.Range("I80").Value = 2509.1365
$ws.Range("J80").Value = 8179
$ws.Range("K80").Value = 7527.4095
$ws.Range("L80").Value = 24537
$ws.Range("M80").Value = -6529.4095
$ws.Range("N80").Value = -26533
$ws.Range("H83").Value = 5136.6343
$ws.Range("I83").Value = 2509.1365
$ws.Range("J83").Value = 8179
$ws.Range("K83").Value = 22582.2285
$ws.Range("L83").Value = 73611
$ws.Range("M83").Value = -17590.2285
$ws.Range("N83").Value = -83595
$ws.Range("H86").Value = 4536.913
$ws.Range("I86").Value = 2355.3333
$ws.Range("J86").Value = 6916.8184
$ws.Range("K86").Value = 2355.3333
$ws.Range("L86").Value = 6916.8184
$ws.Range("M86").Value = -1232.3333
$ws.Range("N86").Value = -9162.8184
$ws.Range("H89").Value = 4536.913
$ws.Range("I89").Value = 2355.3333
$ws.Range("J89").Value = 6916.8184
$ws.Range("K89").Value = 11776.6665
$ws.Range("L89").Value = 34584.092
$ws.Range("M89").Value = -6160.666499999999
$ws.Range("N89").Value = -45816.092
$ws.Range("H116").Value = 4252.9443
$ws.Range("I116").Value = 4735.6
$ws.Range("J116").Value = 3649.625
$ws.Range("K116").Value = 4735.6
$ws.Range("L116").Value = 3649.625
$ws.Range("M116").Value = -1293.6
$ws.Range("N116").Value = -10533.625
$ws.Range("H125").Value = 2816.6667
$ws.Range("J125").Value = 2816.6667
$ws.Range("L125").Value = 25350.0003
$ws.Range("N125").Value = -30270.0003
$ws.Range("H129").Value = 1354.6154
$ws.Range("I129").Value = 277.5
$ws.Range("J129").Value = 1833.3334
$ws.Range("K129").Value = 832.5
$ws.Range("L129").Value = 5500.0002
$ws.Range("M129").Value = 4167.5
$ws.Range("N129").Value = -15500.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1147541.5
$ws.Range("I32").Value = 1416500
$ws.Range("J32").Value = 14073.429
$ws.Range("K32").Value = 1416500
$ws.Range("L32").Value = 14073.429
$ws.Range("M32").Value = -1416213
$ws.Range("N32").Value = -14647.429
$ws.Range("H88").Value = 3890.5454
$ws.Range("I88").Value = 6226.5
$ws.Range("J88").Value = 2555.7144
$ws.Range("K88").Value = 6226.5
$ws.Range("L88").Value = 2555.7144
$ws.Range("M88").Value = -5820.5
$ws.Range("N88").Value = -3367.7144
$ws.Range("H91").Value = 3890.5454
$ws.Range("I91").Value = 6226.5
$ws.Range("J91").Value = 2555.7144
$ws.Range("K91").Value = 6226.5
$ws.Range("L91").Value = 2555.7144
$ws.Range("M91").Value = -4822.5
$ws.Range("N91").Value = -5363.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2980.3225
$ws.Range("I86").Value = 4934.615
$ws.Range("J86").Value = 1568.8889
$ws.Range("K86").Value = 4934.615
$ws.Range("L86").Value = 1568.8889
$ws.Range("M86").Value = -3811.615
$ws.Range("N86").Value = -3814.8889
$ws.Range("H89").Value = 2980.3225
$ws.Range("I89").Value = 4934.615
$ws.Range("J89").Value = 1568.8889
$ws.Range("K89").Value = 24673.075
$ws.Range("L89").Value = 7844.4445
$ws.Range("M89").Value = -19057.075
$ws.Range("N89").Value = -19076.4445
$ws.Range("H99").Value = 10699.833
$ws.Range("I99").Value = 10699.833
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10699.833
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -9201.833000000001
$ws.Range("N99").ClearContents()
$ws.Range("H130").Value = 48672.5
$ws.Range("J130").Value = 48672.5
$ws.Range("L130").Value = 48672.5
$ws.Range("N130").Value = -58712.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3883
$ws.Range("I31").Value = 2958.6667
$ws.Range("J31").Value = 5376.154
$ws.Range("K31").Value = 2958.6667
$ws.Range("L31").Value = 5376.154
$ws.Range("M31").Value = -2663.6667
$ws.Range("N31").Value = -5966.154
$ws.Range("H34").Value = 3883
$ws.Range("I34").Value = 2958.6667
$ws.Range("J34").Value = 5376.154
$ws.Range("K34").Value = 2958.6667
$ws.Range("L34").Value = 5376.154
$ws.Range("M34").Value = -2756.6667
$ws.Range("N34").Value = -5780.154
$ws.Range("H122").Value = 1846.9231
$ws.Range("I122").Value = 1030.2858
$ws.Range("J122").Value = 2799.6667
$ws.Range("K122").Value = 3090.8574
$ws.Range("L122").Value = 8399.000100000001
$ws.Range("M122").Value = -640.8574000000003
$ws.Range("N122").Value = -13299.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 10.612904
$ws.Range("I12").Value = 12.615385
$ws.Range("J12").Value = 9.166667
$ws.Range("K12").Value = 37.846155
$ws.Range("L12").Value = 27.500001
$ws.Range("M12").Value = 135.153845
$ws.Range("N12").Value = -373.500001
$ws.Range("H131").Value = 1386.5892
$ws.Range("I131").Value = 1130.7142
$ws.Range("J131").Value = 1471.881
$ws.Range("K131").Value = 3392.1426
$ws.Range("L131").Value = 4415.643
$ws.Range("M131").Value = 1647.8574
$ws.Range("N131").Value = -14495.643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5241.316
$ws.Range("I70").Value = 4092.6667
$ws.Range("J70").Value = 5456.6875
$ws.Range("K70").Value = 4092.6667
$ws.Range("L70").Value = 5456.6875
$ws.Range("M70").Value = -3822.6667
$ws.Range("N70").Value = -5996.6875
$ws.Range("H73").Value = 5241.316
$ws.Range("I73").Value = 4092.6667
$ws.Range("J73").Value = 5456.6875
$ws.Range("K73").Value = 4092.6667
$ws.Range("L73").Value = 5456.6875
$ws.Range("M73").Value = -3156.6667
$ws.Range("N73").Value = -7328.6875
$ws.Range("H113").Value = 1712.75
$ws.Range("I113").Value = 960.1667
$ws.Range("J113").Value = 2277.1875
$ws.Range("K113").Value = 960.1667
$ws.Range("L113").Value = 2277.1875
$ws.Range("M113").Value = 1209.8333
$ws.Range("N113").Value = -6617.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9590.419
$ws.Range("I132").Value = 3240.3333
$ws.Range("J132").Value = 15543.625
$ws.Range("K132").Value = 9720.999899999999
$ws.Range("L132").Value = 46630.875
$ws.Range("M132").Value = -7190.999899999999
$ws.Range("N132").Value = -51690.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1283.2632
$ws.Range("I126").Value = 848.3570999999999
$ws.Range("J126").Value = 2501
$ws.Range("K126").Value = 2545.0713
$ws.Range("L126").Value = 7503
$ws.Range("M126").Value = -75.07129999999961
$ws.Range("N126").Value = -12443
